$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Insert()

for ($i = 2; $i -le 30; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 1).Font.Bold = $ws.Cells.Item($i, 2).Font.Bold
    $ws.Cells.Item($i, 1).HorizontalAlignment = $ws.Cells.Item($i, 2).HorizontalAlignment
    $ws.Cells.Item($i, 1).VerticalAlignment = $ws.Cells.Item($i, 2).VerticalAlignment
    $ws.Cells.Item($i, 1).ShrinkToFit = $ws.Cells.Item($i, 2).ShrinkToFit
    $ws.Cells.Item($i, 1).Borders.LineStyle = $ws.Cells.Item($i, 2).Borders.LineStyle
}
